$wb = $excel.ActiveWorkbook

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A50").Value = "'2980"
    $ws.Range("A50").Style = "Normal"

    $ws.Range("B50").Value = "'2025-09-10"
    $ws.Range("B50").Style = "Normal"

    $ws.Range("C50").Value = "Erdemli"

    $ws.Range("D50").Value = "'1"
    $ws.Range("D50").Style = "Normal"

    $ws.Range("E50").Value = "MAKS"

    $ws.Range("F50").Value = "SERDAR ARSLAN (Tekniker), ÖZKAN AKBAŞ (Mühendis)"
}
